$d = $word.ActiveDocument
$failCount = 0

# Paragraph 6, Run 1
$rng = $d.Paragraphs.Item(6).Range
$found = $rng.Find.Execute("Apresentar aos alunos a Engenharia Bioquímica, as características da profissão e orientar quanto as atribuições e as áreas de atuação do Engenheiro Bioquímico. Além disso, desenvolver nos alunos uma visão macro dos tipos e etapas de um bioprocesso industrial e, por fim, orientar sobre a atuação do Engenheiro Bioquímico na indústria, pesquisa e ensino, e empreendedorismo e inovação em engenharia.", $true, $false, $false, $false, $false, $true, 0, $false, "1. Histórico da Engenharia Bioquímica 2. Engenharia Bioquímica: Definições e conceitos 3. Mercado de trabalho de Engenharia 4. Áreas de atuação do Engenheiro Bioquímico 5. A Indústria de Bioprocessos 6. Escalas de produção 7. Estudo de casos (processos biotecnológicos) 8. Visita supervisionada.", 2)
if (-not $found) { Write-Output "FAILED: P6R1"; $failCount++ }

# Paragraph 7, Run 1
$rng = $d.Paragraphs.Item(7).Range
$found = $rng.Find.Execute("Introduce students to Biochemical Engineering, the characteristics of the profession, and guide them regarding the responsibilities and areas of practice of a Biochemical Engineer. Additionally, develop in students a macro view of the types and stages of an industrial bioprocess, and finally, guide them on the role of the Biochemical Engineer in industry, research and teaching, as well as entrepreneurship and innovation in engineering.", $true, $false, $false, $false, $false, $true, 0, $false, "1. History of Biochemical Engineering; 2. Biochemical Engineering: Definitions and concepts; 3. Job market for Engineering; 4. Areas of practice for the Biochemical Engineer; 5. The Bioprocess Industry; 6. Production scales^lCase studies (biotechnological processes); 7. Supervised visit.", 2)
if (-not $found) { Write-Output "FAILED: P7R1"; $failCount++ }

# Paragraph 9, Run 2
$rng = $d.Paragraphs.Item(9).Range
$found = $rng.Find.Execute("5817181 - Valdeir Arantes", $true, $false, $false, $false, $false, $true, 0, $false, "1. Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2. Mercado de trabalho da Engenharia do Brasil 3. Atribuições e áreas de atuação do Engenheiro Bioquímico 4. Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5. A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6. Escalas de produção – laboratório, piloto, industrial. 7. Estudo de casos (processos biotecnológicos). 8. Empreendedorismo e Inovação em Engenharia. 9. Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso.1. Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2. Mercado de trabalho da Engenharia do Brasil 3. Atribuições e áreas de atuação do Engenheiro Bioquímico 4. Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5. A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6. Escalas de produção – laboratório, piloto, industrial. 7. Estudo de casos (processos biotecnológicos). 8. Empreendedorismo e Inovação em Engenharia. 9. Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso.", 2)
if (-not $found) { Write-Output "FAILED: P9R2"; $failCount++ }

# Paragraph 9, Run 1
$rng = $d.Paragraphs.Item(9).Range
$found = $rng.Find.Execute("101761 - Arnaldo Márcio Ramalho Prata^l", $true, $false, $false, $false, $false, $true, 0, $false, "Apresentar aos alunos a Engenharia Bioquímica, as características da profissão e orientar quanto as atribuições e as áreas de atuação do Engenheiro Bioquímico. Além disso, desenvolver nos alunos uma visão macro dos tipos e etapas de um bioprocesso industrial e, por fim, orientar sobre a atuação do Engenheiro Bioquímico na indústria, pesquisa e ensino, e empreendedorismo e inovação em engenharia.^l", 2)
if (-not $found) { Write-Output "FAILED: P9R1"; $failCount++ }

# Paragraph 11, Run 1
$rng = $d.Paragraphs.Item(11).Range
$found = $rng.Find.Execute("1. Histórico da Engenharia Bioquímica 2. Engenharia Bioquímica: Definições e conceitos 3. Mercado de trabalho de Engenharia 4. Áreas de atuação do Engenheiro Bioquímico 5. A Indústria de Bioprocessos 6. Escalas de produção 7. Estudo de casos (processos biotecnológicos) 8. Visita supervisionada.", $true, $false, $false, $false, $false, $true, 0, $false, "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina; exercícios; dinâmicas. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a aplicações dos conceitos abordados à um processo, produto ou serviço na área de Engenharia de Bioquímica e que relacione com a formação acadêmica e atribuições profissionais do Engenheiro Bioquímico.", 2)
if (-not $found) { Write-Output "FAILED: P11R1"; $failCount++ }

# Paragraph 12, Run 1
$rng = $d.Paragraphs.Item(12).Range
$found = $rng.Find.Execute("1. History of Biochemical Engineering; 2. Biochemical Engineering: Definitions and concepts; 3. Job market for Engineering; 4. Areas of practice for the Biochemical Engineer; 5. The Bioprocess Industry; 6. Production scales^lCase studies (biotechnological processes); 7. Supervised visit.", $true, $false, $false, $false, $false, $true, 0, $false, "Introduce students to Biochemical Engineering, the characteristics of the profession, and guide them regarding the responsibilities and areas of practice of a Biochemical Engineer. Additionally, develop in students a macro view of the types and stages of an industrial bioprocess, and finally, guide them on the role of the Biochemical Engineer in industry, research and teaching, as well as entrepreneurship and innovation in engineering.", 2)
if (-not $found) { Write-Output "FAILED: P12R1"; $failCount++ }

# Paragraph 14, Run 1
$rng = $d.Paragraphs.Item(14).Range
$found = $rng.Find.Execute("1. Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2. Mercado de trabalho da Engenharia do Brasil 3. Atribuições e áreas de atuação do Engenheiro Bioquímico 4. Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5. A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6. Escalas de produção – laboratório, piloto, industrial. 7. Estudo de casos (processos biotecnológicos). 8. Empreendedorismo e Inovação em Engenharia. 9. Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso.1. Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2. Mercado de trabalho da Engenharia do Brasil 3. Atribuições e áreas de atuação do Engenheiro Bioquímico 4. Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5. A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6. Escalas de produção – laboratório, piloto, industrial. 7. Estudo de casos (processos biotecnológicos). 8. Empreendedorismo e Inovação em Engenharia. 9. Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso.", $true, $false, $false, $false, $false, $true, 0, $false, "A nota (N) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.", 2)
if (-not $found) { Write-Output "FAILED: P14R1"; $failCount++ }

# Paragraph 17, Run 6
$rng = $d.Paragraphs.Item(17).Range
$found = $rng.Find.Execute("Média Final = (N + Prova Recuperação)/2", $true, $false, $false, $false, $false, $true, 0, $false, "101761 - Arnaldo Márcio Ramalho Prata", 2)
if (-not $found) { Write-Output "FAILED: P17R6"; $failCount++ }

# Paragraph 17, Run 4
$rng = $d.Paragraphs.Item(17).Range
$found = $rng.Find.Execute("A nota (N) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.^l", $true, $false, $false, $false, $false, $true, 0, $false, "Schmidell, W.; Lima, U. A.; Aquarone, E.; Borzani, W. Biotecnologia Industrial – Engenharia^lBioquímica, vol. 2, São Paulo: Edgard Blücher, 2001.^lShuler, L. M.; Kargi, F. Bioprocess Engineering – Basic Concepts. Second edition. New^lJersey: PrenticeHall,^l2002.^lArigos atuais relacionaos com o tema de Engenharia Bioquímica^l", 2)
if (-not $found) { Write-Output "FAILED: P17R4"; $failCount++ }

# Paragraph 17, Run 2
$rng = $d.Paragraphs.Item(17).Range
$found = $rng.Find.Execute("O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina; exercícios; dinâmicas. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a aplicações dos conceitos abordados à um processo, produto ou serviço na área de Engenharia de Bioquímica e que relacione com a formação acadêmica e atribuições profissionais do Engenheiro Bioquímico.^l", $true, $false, $false, $false, $false, $true, 0, $false, "Média Final = (N + Prova Recuperação)/2^l", 2)
if (-not $found) { Write-Output "FAILED: P17R2"; $failCount++ }

# Paragraph 19, Run 1
$rng = $d.Paragraphs.Item(19).Range
$found = $rng.Find.Execute("Schmidell, W.; Lima, U. A.; Aquarone, E.; Borzani, W. Biotecnologia Industrial – Engenharia^lBioquímica, vol. 2, São Paulo: Edgard Blücher, 2001.^lShuler, L. M.; Kargi, F. Bioprocess Engineering – Basic Concepts. Second edition. New^lJersey: PrenticeHall,^l2002.^lArigos atuais relacionaos com o tema de Engenharia Bioquímica", $true, $false, $false, $false, $false, $true, 0, $false, "5817181 - Valdeir Arantes", 2)
if (-not $found) { Write-Output "FAILED: P19R1"; $failCount++ }

Write-Output "Done. Failures: $failCount"